$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = '$2b$12$8/Iyol5kVAkJijBh9CVMhuEeMlNKHBKZAfa1GvcMtyFcD1BeuDIfy'
$ws.Range("B3").Value = '$2b$12$XFVj25dIinW2VVzVjzM5A.r9CHfH2opXiSyZw1wcP0naySiOJQeHO'
$ws.Range("B4").Value = '$2b$12$dNuOuP8WabQtRuuba8sadeOD9NsRpZBu3zHSpG5a/d04ZIOO3OLMa'
$ws.Range("B5").Value = '$2b$12$B40wy1tIkOaJ7ZQDNAFG0O0j6gTFxowXbOi2ct.PUdtpqbH4Y/2Ui'
